$wb = $excel.ActiveWorkbook

$wsStd = $wb.Worksheets.Item("Stundenerfassung")
$wsProj = $wb.Worksheets.Item("Projekte")

# --- Append the two new "Stundenerfassung" rows (115 & 116) ---------------
# Row 115: 21.08.2017 (serial 42968), ETIC2 / Codierung nach MVVM, 3h
$wsStd.Range("A115").Value = 42968
$wsStd.Range("A114").Copy() | Out-Null
$wsStd.Range("A115").PasteSpecial(-4122) | Out-Null   # xlPasteFormats (keep date format s="3")

$wsStd.Range("B115").Value = "ETIC2"
$wsStd.Range("C115").Value = "Codierung nach MVVM"
$wsStd.Range("D115").Value = 3

# Row 116: 21.08.2017 (serial 42968), Abgabe / Diverses, 4h
$wsStd.Range("A116").Value = 42968
$wsStd.Range("A114").Copy() | Out-Null
$wsStd.Range("A116").PasteSpecial(-4122) | Out-Null   # xlPasteFormats (keep date format s="3")

$wsStd.Range("B116").Value = "Abgabe"
$wsStd.Range("C116").Value = "Diverses"
$wsStd.Range("D116").Value = 4

# --- View/selection bookkeeping --------------------------------------------
# Projekte keeps its own scroll position but loses the tab-selection and
# moves its cell-selection to D17 (do this first, it's not the final tab)
$wsProj.Activate()
$wsProj.Range("D17").Select() | Out-Null

# Stundenerfassung becomes the active sheet/tab again (was Projekte) - do
# this last so it ends up as the saved, active tab
$wsStd.Activate()
$aw = $excel.ActiveWindow
$aw.ScrollRow = 105
$aw.ScrollColumn = 1
$wsStd.Range("E115").Select() | Out-Null
